# step 4: final changes of the day
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Highlight row 2 (A2:F2) with a yellow fill
$ws.Range("A2:F2").Interior.Color = 65535

# Make row 3 (A3:F3) bold
$ws.Range("A3:F3").Font.Bold = $true

# Update the labels in column F to describe the new formatting being demoed
$ws.Range("F2").Value = "highlight"
$ws.Range("F3").Value = "bold"
$ws.Range("F4").Value = "conditional formatting"

# Add a conditional formatting rule to A4:F4 - highlight (red) when A4 > 0
$rng = $ws.Range("A4:F4")
$fc = $rng.FormatConditions.Add(2, 0, "=`$A`$4>0")
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615

# Resize column F to fit its new, longer content
$ws.Columns("F").AutoFit()

# Leave the A4:F4 row selected, as that's what was being worked on
$ws.Range("A4:F4").Select()
